{"js": "// Remove specific course-requirement lines from the \"Requisitos\" list.\n// Each line in that list is a run ending in a manual line break\n// (U+000B, Word's \"\\v\") inside a single paragraph, so we search for the\n// exact line text including its trailing break and delete the match \u2014\n// this removes both the text and the line break that followed it.\nconst coursesToRemove = [\n  \"LOM3081 -  Introdu\u00e7\u00e3o \u00e0 Mec\u00e2nica dos S\u00f3lidos  (Requisito)\",\n  \"LOM3205 -  Eletromagnetismo  (Requisito)\",\n  \"LOM3212 -  Fen\u00f4menos de Transporte A  (Requisito)\",\n  \"LOM3240 -  Qu\u00edmica Inorg\u00e2nica Fundamental e Aplicada  (Requisito)\",\n  \"LOM3253 -  F\u00edsica Matem\u00e1tica  (Requisito)\",\n  \"LOM3257 -  Mec\u00e2nica Cl\u00e1ssica  (Requisito)\",\n  \"LOM3262 -  Circuitos El\u00e9tricos  (Requisito)\"\n];\n\nconst body = context.document.body;\n\nfor (const course of coursesToRemove) {\n  const results = body.search(course + \"\\u000b\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].delete();\n  }\n  await context.sync();\n}\n", "ps1": "# Remove specific course-requirement lines from the \"Requisitos\" list.\n# Each line in that list is terminated by a manual line break (vertical\n# tab / Word \"Line\" break) rather than a paragraph mark, so we locate each\n# line's text with Find, extend the found range by one character to grab\n# the trailing line break, and delete the whole range - this removes both\n# the text and the break that followed it, leaving the rest of the list\n# intact.\n$d = $word.ActiveDocument\n\n$coursesToRemove = @(\n  \"LOM3081 -  Introdu\u00e7\u00e3o \u00e0 Mec\u00e2nica dos S\u00f3lidos  (Requisito)\",\n  \"LOM3205 -  Eletromagnetismo  (Requisito)\",\n  \"LOM3212 -  Fen\u00f4menos de Transporte A  (Requisito)\",\n  \"LOM3240 -  Qu\u00edmica Inorg\u00e2nica Fundamental e Aplicada  (Requisito)\",\n  \"LOM3253 -  F\u00edsica Matem\u00e1tica  (Requisito)\",\n  \"LOM3257 -  Mec\u00e2nica Cl\u00e1ssica  (Requisito)\",\n  \"LOM3262 -  Circuitos El\u00e9tricos  (Requisito)\"\n)\n\nforeach ($course in $coursesToRemove) {\n  $rng = $d.Content\n  $rng.Find.MatchCase = $true\n  $found = $rng.Find.Execute($course)\n  if ($found) {\n    $rng.MoveEnd(1, 1)\n    $rng.Delete()\n  }\n}\n"}
